$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.291.65"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "1.721.70"
$ws.Range("E3").Value = "  +3.33%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("D5").Value = "'240.32"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4736"
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("D8").Value = "'0.2633"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'0.06198"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "1.717.88"
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("D11").Value = "'0.07058"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("E12").Value = "  +4.84%  "
$ws.Range("D13").Value = "'0.5996"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "'4.440"
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "'76.29"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "'0.9996"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "26.294.49"
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("D19").Value = "'0.000006818"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D21").Value = "1.937.23"
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("D22").Value = "'4.537"
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("D23").Value = "'8.729"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'5.255"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "'135.15"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "'15.22"
$ws.Range("D27").Value = "'1.769"
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("D28").Value = "'1.400"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").Value = "'106.73"
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("D30").Value = "'3.956"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "'3.689"
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("D32").Value = "'0.07803"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "'0.04512"
$ws.Range("E33").Value = "  +7.07%  "
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "'0.9992"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.614"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9829"
$ws.Range("E36").Value = "  +3.81%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.6243"
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'0.9336"
$ws.Range("E38").Value = "  +8.32%  "
$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").Value = "'114.86"
$ws.Range("E39").Value = "  +18.17%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.455"
$ws.Range("E40").Value = "  -5.44%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'1.930"
$ws.Range("E41").Value = "  +4.45%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'0.9999"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.657"
$ws.Range("E43").Value = "  +16.86%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.01485"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.3841"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1186"
$ws.Range("E46").Value = "  +5.95%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'6.365"
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05269"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.898"
$ws.Range("E49").Value = "  +7.85%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'30.42"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3390"
